$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.02539962882261295; C = 0.2597273962171842; D = 0.3180567948316633; E = 0.5639652425740999; F = 0.5688896304229624 }
    3  = @{ B = 0.4459617371143483;  C = 0.6540342364242132; D = 1.388895812590984;  E = 1.17851423945194;   F = 1.101732372426877 }
    4  = @{ B = 0.6024413171113788;  C = 0.9796726598280953; D = 3.63049943662146;   E = 1.905386951939543;  F = 1.82599257214512 }
    5  = @{ B = 0.5121426763354459;  C = 1.006154916801229;  D = 3.970987900588087;  E = 1.99273377564292;   F = 1.945754947762022;  G = 49 }
    6  = @{ B = 0.4218352939147489;  C = 0.8877411355809319; D = 3.626482144461433;  E = 1.904332466892646;  F = 1.876675320794537;  G = 48 }
    7  = @{ B = 0.3949550873240204;  C = 0.8680554327475212; D = 4.006566611836582;  E = 2.001640979755506;  F = 1.987940659805524;  G = 39 }
    8  = @{ B = 0.3899578886789405;  C = 0.8693687210883813; D = 4.090339984385427;  E = 2.022458895598481;  F = 2.011147094482062;  G = 38 }
    9  = @{ B = 0.3439843485875357;  C = 1.183619751688364;  D = 6.937012254809386;  E = 2.633820847136226;  F = 2.675746881502995;  G = 21 }
    10 = @{ B = -0.2090860748678977; C = 0.8864813411966667; D = 3.164946506954499;  E = 1.779029653197074;  F = 1.833391419946016;  G = 14 }
    11 = @{ B = 0.528408343223261;   C = 0.528408343223261;  D = 0.3628698470762892; E = 0.6023867919171944; F = 0.3233698924767455 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
